$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string]) {
            if ($val -match "Patrick Toole") {
                $newVal = $val -replace "Patrick Toole III", "Patrick III Toole"
                $newVal = $newVal -replace "Patrick Toole II", "Patrick II Toole"
                $cell.Value2 = $newVal
            }
        }
    }
}

$ws.Range("E26").Select() | Out-Null
